# Actualización 11 de Mayo - Mañana
#
# Refresh of the "Rescatables" (make-up/rescue candidates) sheet: the
# roster gains a new student (TREJO LUENGAS ELIZABETH, mat. 19330051920414,
# group 4ASV, 2 reprobadas) and the previously-listed students are
# re-sorted/updated to their latest values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Insert a new row at position 4, shifting the remaining students down by
# one. This grows the table from 7 data rows (rows 2-7) to 8 data rows
# (rows 2-8), leaving room for the new student.
$ws.Rows.Item(4).Insert()

# Final values for the data block (rows 2-8), column by column.
$mat    = @(19330051920366, 19330051920266, 19330051920414, 19330051920287, 19330051920330, 19330051920351, 19330051920372)
$pat    = @("CRISTOBAL", "MIXCOHUA", "TREJO", "MARQUEZ", "LOPEZ", "TRUJILLO", "FLORES")
$mater  = @("BRUNO", "IXMATLAHUA", "LUENGAS", "DE JESUS", "TZOPITL", "DE LA LUZ", "PEREZ")
$nom    = @("DANIELA", "VALENTIN", "ELIZABETH", "ANGEL JARET", "DIEGO", "EDGAR", "ANTONIO")
$largo  = @("CÁLCULO DIFERENCIAL", "CÁLCULO DIFERENCIAL", "CÁLCULO DIFERENCIAL", "CÁLCULO DIFERENCIAL", "CÁLCULO DIFERENCIAL", "CÁLCULO DIFERENCIAL", "CÁLCULO DIFERENCIAL")
$grupo  = @("4ARHV", "4ASV", "4ASV", "4ALCV", "4APV", "4APV", "4ARHV")
$reprob = @(2, 2, 2, 1, 1, 1, 1)

for ($i = 0; $i -lt $mat.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $mat[$i]
    $ws.Cells.Item($row, 2).Value = $pat[$i]
    $ws.Cells.Item($row, 3).Value = $mater[$i]
    $ws.Cells.Item($row, 4).Value = $nom[$i]
    $ws.Cells.Item($row, 5).Value = $largo[$i]
    $ws.Cells.Item($row, 6).Value = $grupo[$i]
    $ws.Cells.Item($row, 7).Value = $reprob[$i]
}
